$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 19
$ws.Range("AO4").Value = 10
$ws.Range("G4").Value = 1.8
$ws.Range("AI5").Value = 15
$ws.Range("AJ5").Value = 13
$ws.Range("G5").Value = 2.5
$ws.Range("I5").Value = 3.2
$ws.Range("J5").Value = 3.25
$ws.Range("Z5").Value = 23
$ws.Range("Q7").Value = 1.98
$ws.Range("R7").Value = 1.88
$ws.Range("AA8").Value = 34
$ws.Range("AH8").Value = 7.5
$ws.Range("AI8").Value = 8.5
$ws.Range("AK8").Value = 15
$ws.Range("AN8").Value = 6
$ws.Range("AP8").Value = 29
$ws.Range("AT8").Value = 2.75
$ws.Range("AX8").Value = 9.5
$ws.Range("G8").Value = 4.2
$ws.Range("I8").Value = 1.73
$ws.Range("J8").Value = 4.75
$ws.Range("L8").Value = 2.4
$ws.Range("M8").Value = 1.05
$ws.Range("O8").Value = 1.29
$ws.Range("S8").Value = 1.4
$ws.Range("T8").Value = 2.75
$ws.Range("U9").Value = 1.63
$ws.Range("U11").Value = 1.58
$ws.Range("AF12").Value = 41
$ws.Range("AH12").Value = 13
$ws.Range("AJ12").Value = 12
$ws.Range("AK12").Value = 34
$ws.Range("AL12").Value = 23
$ws.Range("AM12").Value = 29
$ws.Range("AT12").Value = 3.25
$ws.Range("AU12").Value = 7.5
$ws.Range("G12").Value = 2.2
$ws.Range("I12").Value = 3.1
$ws.Range("J12").Value = 2.75
$ws.Range("L12").Value = 3.5
$ws.Range("S12").Value = 1.33
$ws.Range("T12").Value = 3.25
$ws.Range("U12").Value = 1.57
$ws.Range("V12").Value = 2.25
$ws.Range("X12").Value = 12
$ws.Range("Y12").Value = 9
$ws.Range("Z12").Value = 21
$ws.Range("M13").Value = 1.02
$ws.Range("N13").Value = 21
$ws.Range("I14").Value = 1.7
$ws.Range("G15").Value = 1.75
$ws.Range("Q15").Value = 1.67
$ws.Range("R15").Value = 2.15
$ws.Range("G16").Value = 1.53
$ws.Range("Q19").Value = 1.67
$ws.Range("Q20").Value = 1.48
$ws.Range("Q21").Value = 1.33
$ws.Range("Q22").Value = 1.7
$ws.Range("M23").Value = 1.05
$ws.Range("O23").Value = 1.41
$ws.Range("P23").Value = 2.62
$ws.Range("R23").Value = 1.57
$ws.Range("M24").Value = 1.05
$ws.Range("O24").Value = 1.37
$ws.Range("R24").Value = 1.6
$ws.Range("G25").Value = 1.83
$ws.Range("M25").Value = 1.03
$ws.Range("O25").Value = 1.22
$ws.Range("M26").Value = 1.03
$ws.Range("O26").Value = 1.19
$ws.Range("Q26").Value = 1.8
$ws.Range("AA27").Value = 19
$ws.Range("AC27").Value = 7
$ws.Range("AF27").Value = 67
$ws.Range("AI27").Value = 21
$ws.Range("AL27").Value = 41
$ws.Range("AP27").Value = 26
$ws.Range("AU27").Value = 9
$ws.Range("BA27").Value = 126
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 3
$ws.Range("I27").Value = 4.1
$ws.Range("K27").Value = 2
$ws.Range("M27").Value = 1.1
$ws.Range("N27").Value = 7
$ws.Range("O27").Value = 1.4
$ws.Range("P27").Value = 2.75
$ws.Range("Q27").Value = 2.35
$ws.Range("R27").Value = 1.57
$ws.Range("V27").Value = 1.73
$ws.Range("W27").Value = 6
$ws.Range("Y27").Value = 9.5
$ws.Range("U28").Value = 1.8
$ws.Range("V28").Value = 1.91
$ws.Range("M30").Value = 1.05
$ws.Range("O30").Value = 1.27
$ws.Range("U30").Value = 1.8
$ws.Range("V30").Value = 1.91
$ws.Range("M32").Value = 1.02
$ws.Range("O32").Value = 1.13
$ws.Range("AE33").Value = 15
$ws.Range("AG33").Value = 351
$ws.Range("AH33").Value = 7
$ws.Range("AI33").Value = 11
$ws.Range("AJ33").Value = 10
$ws.Range("AK33").Value = 23
$ws.Range("AQ33").Value = 51
$ws.Range("AR33").Value = 81
$ws.Range("BA33").Value = 67
$ws.Range("BD33").Value = 126
$ws.Range("G33").Value = 3.2
$ws.Range("I33").Value = 2.35
$ws.Range("J33").Value = 3.75
$ws.Range("M33").Value = 1.05
$ws.Range("O33").Value = 1.37
$ws.Range("N34").Value = 12
$ws.Range("U34").Value = 1.87
$ws.Range("V34").Value = 1.87
$ws.Range("AG35").Value = 301
$ws.Range("AO35").Value = 15
$ws.Range("AQ35").Value = 51
$ws.Range("AW35").Value = 4.75
$ws.Range("G35").Value = 2.45
$ws.Range("I35").Value = 3.1
$ws.Range("J35").Value = 3.2
$ws.Range("L35").Value = 3.6
$ws.Range("M35").Value = 1.08
$ws.Range("N35").Value = 8
$ws.Range("O35").Value = 1.36
$ws.Range("U35").Value = 1.8
$ws.Range("V35").Value = 1.8
$ws.Range("W35").Value = 7.5
$ws.Range("M36").Value = 1.05
$ws.Range("O36").Value = 1.29
$ws.Range("M38").Value = 1.08
$ws.Range("O38").Value = 1.44
$ws.Range("P38").Value = 2.63
$ws.Range("M39").Value = 1.05
$ws.Range("O39").Value = 1.29
